$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 25,9
$arr[0,0] = "model_7_5_24"
$arr[0,1] = 0.7247495437567779
$arr[0,2] = 0.06376739587406144
$arr[0,3] = 0.82409504076811
$arr[0,4] = 0.3847048822927533
$arr[0,5] = 0.304621160030365
$arr[0,6] = 1.549709916114807
$arr[0,7] = 0.2319634109735489
$arr[0,8] = 0.9295945763587952
$arr[1,0] = "model_7_5_23"
$arr[1,1] = 0.7283151784190296
$arr[1,2] = 0.07765003039167628
$arr[1,3] = 0.8302995674595636
$arr[1,4] = 0.3953069529429234
$arr[1,5] = 0.3006750643253326
$arr[1,6] = 1.526730537414551
$arr[1,7] = 0.2237815707921982
$arr[1,8] = 0.9135769605636597
$arr[2,0] = "model_7_5_22"
$arr[2,1] = 0.7322299905351446
$arr[2,2] = 0.09324926030642955
$arr[2,3] = 0.8363128304565055
$arr[2,4] = 0.4068222504097677
$arr[2,5] = 0.2963425517082214
$arr[2,6] = 1.500909805297852
$arr[2,7] = 0.2158519625663757
$arr[2,8] = 0.8961794972419739
$arr[3,0] = "model_7_5_21"
$arr[3,1] = 0.7366564787997097
$arr[3,2] = 0.1110998820209144
$arr[3,3] = 0.8423777192207784
$arr[3,4] = 0.4196689343026407
$arr[3,5] = 0.2914436757564545
$arr[3,6] = 1.471362352371216
$arr[3,7] = 0.207854300737381
$arr[3,8] = 0.876770555973053
$arr[4,0] = "model_7_5_20"
$arr[4,1] = 0.7415264976253521
$arr[4,2] = 0.1309114170639998
$arr[4,3] = 0.848536849800671
$arr[4,4] = 0.4336901501606907
$arr[4,5] = 0.2860539853572845
$arr[4,6] = 1.438569068908691
$arr[4,7] = 0.1997323334217072
$arr[4,8] = 0.8555871844291687
$arr[5,0] = "model_7_5_19"
$arr[5,1] = 0.7468641642014335
$arr[5,2] = 0.152921345870812
$arr[5,3] = 0.8545328319623225
$arr[5,4] = 0.4489203907476186
$arr[5,5] = 0.2801467776298523
$arr[5,6] = 1.40213680267334
$arr[5,7] = 0.1918255090713501
$arr[5,8] = 0.8325772285461426
$arr[6,0] = "model_7_5_18"
$arr[6,1] = 0.752883638722313
$arr[6,2] = 0.1777542516197554
$arr[6,3] = 0.8607876574019349
$arr[6,4] = 0.4658921829782201
$arr[6,5] = 0.2734850347042084
$arr[6,6] = 1.361032009124756
$arr[6,7] = 0.1835773587226868
$arr[6,8] = 0.8069360256195068
$arr[7,0] = "model_7_5_17"
$arr[7,1] = 0.7593863819171758
$arr[7,2] = 0.2049069878355617
$arr[7,3] = 0.8668175930059645
$arr[7,4] = 0.4841174374567468
$arr[7,5] = 0.2662883698940277
$arr[7,6] = 1.316087007522583
$arr[7,7] = 0.1756257712841034
$arr[7,8] = 0.7794010639190674
$arr[8,0] = "model_7_5_16"
$arr[8,1] = 0.7663057795530505
$arr[8,2] = 0.2338371309332362
$arr[8,3] = 0.8730547286292003
$arr[8,4] = 0.5034597314993571
$arr[8,5] = 0.2586306035518646
$arr[8,6] = 1.268200039863586
$arr[8,7] = 0.1674009561538696
$arr[8,8] = 0.7501785755157471
$arr[9,0] = "model_7_5_15"
$arr[9,1] = 0.774408358093261
$arr[9,2] = 0.2676764579664317
$arr[9,3] = 0.879557137700689
$arr[9,4] = 0.5257597832123682
$arr[9,5] = 0.249663457274437
$arr[9,6] = 1.212187170982361
$arr[9,7] = 0.1588263213634491
$arr[9,8] = 0.7164874076843262
$arr[10,0] = "model_7_5_14"
$arr[10,1] = 0.7834092300040275
$arr[10,2] = 0.3053892226534659
$arr[10,3] = 0.8862110057372222
$arr[10,4] = 0.5503666935135632
$arr[10,5] = 0.2397021502256393
$arr[10,6] = 1.149762630462646
$arr[10,7] = 0.1500519514083862
$arr[10,8] = 0.6793109774589539
$arr[11,0] = "model_7_5_13"
$arr[11,1] = 0.7934246207252602
$arr[11,2] = 0.3473317888983217
$arr[11,3] = 0.8930790166467704
$arr[11,4] = 0.5775134858641009
$arr[11,5] = 0.2286181002855301
$arr[11,6] = 1.080336689949036
$arr[11,7] = 0.1409952044487
$arr[11,8] = 0.6382972598075867
$arr[12,0] = "model_7_5_12"
$arr[12,1] = 0.8050318982823272
$arr[12,2] = 0.3960310186087773
$arr[12,3] = 0.9001962825826963
$arr[12,4] = 0.6086835686852032
$arr[12,5] = 0.2157722264528275
$arr[12,6] = 0.9997267127037048
$arr[12,7] = 0.1316097676753998
$arr[12,8] = 0.5912051796913147
$arr[13,0] = "model_7_5_11"
$arr[13,1] = 0.817053324203673
$arr[13,2] = 0.4466141876580161
$arr[13,3] = 0.907202540371934
$arr[13,4] = 0.6409062929581394
$arr[13,5] = 0.202468067407608
$arr[13,6] = 0.9159983396530151
$arr[13,7] = 0.1223707050085068
$arr[13,8] = 0.5425227880477905
$arr[14,0] = "model_7_5_10"
$arr[14,1] = 0.8290776653281506
$arr[14,2] = 0.4975465033471456
$arr[14,3] = 0.9140916799569087
$arr[14,4] = 0.6732755162018431
$arr[14,5] = 0.1891606748104095
$arr[14,6] = 0.8316920399665833
$arr[14,7] = 0.113286092877388
$arr[14,8] = 0.4936189651489258
$arr[15,0] = "model_7_5_9"
$arr[15,1] = 0.841610869147156
$arr[15,2] = 0.5510797141949351
$arr[15,3] = 0.9207342671646286
$arr[15,4] = 0.7070542227921628
$arr[15,5] = 0.1752900928258896
$arr[15,6] = 0.7430805563926697
$arr[15,7] = 0.1045266091823578
$arr[15,8] = 0.4425857365131378
$arr[16,0] = "model_7_5_8"
$arr[16,1] = 0.8538103386396062
$arr[16,2] = 0.603557863011282
$arr[16,3] = 0.9270979761446071
$arr[16,4] = 0.740107451331429
$arr[16,5] = 0.1617888659238815
$arr[16,6] = 0.6562154293060303
$arr[16,7] = 0.09613487869501114
$arr[16,8] = 0.3926485478878021
$arr[17,0] = "model_7_5_7"
$arr[17,1] = 0.8659047125460969
$arr[17,2] = 0.6561778840938759
$arr[17,3] = 0.9331042990212493
$arr[17,4] = 0.7730946982284916
$arr[17,5] = 0.1484039574861526
$arr[17,6] = 0.5691155195236206
$arr[17,7] = 0.08821442723274231
$arr[17,8] = 0.3428110480308533
$arr[18,0] = "model_7_5_6"
$arr[18,1] = 0.8769125906037676
$arr[18,2] = 0.704841064951601
$arr[18,3] = 0.9386698864674884
$arr[18,4] = 0.8036083515405559
$arr[18,5] = 0.1362214982509613
$arr[18,6] = 0.4885652959346771
$arr[18,7] = 0.08087515830993652
$arr[18,8] = 0.2967106699943542
$arr[19,0] = "model_7_5_5"
$arr[19,1] = 0.8868052547222347
$arr[19,2] = 0.749617847544491
$arr[19,3] = 0.9435670692414279
$arr[19,4] = 0.8315926039728895
$arr[19,5] = 0.1252732276916504
$arr[19,6] = 0.4144479930400848
$arr[19,7] = 0.07441731542348862
$arr[19,8] = 0.2544317543506622
$arr[20,0] = "model_7_5_4"
$arr[20,1] = 0.8954115130359251
$arr[20,2] = 0.7903815812439933
$arr[20,3] = 0.9470872057996911
$arr[20,4] = 0.8566818579159845
$arr[20,5] = 0.1157486438751221
$arr[20,6] = 0.346973329782486
$arr[20,7] = 0.06977535784244537
$arr[20,8] = 0.2165266275405884
$arr[21,0] = "model_7_5_3"
$arr[21,1] = 0.9004285761578044
$arr[21,2] = 0.8164053135431721
$arr[21,3] = 0.950424639442038
$arr[21,4] = 0.8731469648246083
$arr[21,5] = 0.1101962253451347
$arr[21,6] = 0.3038972616195679
$arr[21,7] = 0.06537432968616486
$arr[21,8] = 0.1916509866714478
$arr[22,0] = "model_7_5_0"
$arr[22,1] = 0.9032002323269125
$arr[22,2] = 0.8557186389250049
$arr[22,3] = 0.9463487679164062
$arr[22,4] = 0.8942760518742964
$arr[22,5] = 0.1071288213133812
$arr[22,6] = 0.238823413848877
$arr[22,7] = 0.07074913382530212
$arr[22,8] = 0.1597288995981216
$arr[23,0] = "model_7_5_2"
$arr[23,1] = 0.9037192653810219
$arr[23,2] = 0.8372102935678493
$arr[23,3] = 0.9515494534337017
$arr[23,4] = 0.8856763609933811
$arr[23,5] = 0.1065543964505196
$arr[23,6] = 0.2694595754146576
$arr[23,7] = 0.06389105319976807
$arr[23,8] = 0.1727214455604553
$arr[24,0] = "model_7_5_1"
$arr[24,1] = 0.9048022956938477
$arr[24,2] = 0.8510058844790237
$arr[24,3] = 0.9501462446213255
$arr[24,4] = 0.8931018893132748
$arr[24,5] = 0.1053558066487312
$arr[24,6] = 0.2466242611408234
$arr[24,7] = 0.06574144959449768
$arr[24,8] = 0.1615028530359268

$ws.Range("A2:I26").Value = $arr
